$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update productname on both sheets (B1)
$ws1.Range("B1").Value = "2500-RBI-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-UPF-1st"
$ws2.Range("B1").Value = "2500-RBI-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-2-DATE-VAR-INST-UPF-1st"

# Update shortname on ProductLoanInput (B2) - now text instead of numeric
$ws1.Range("B2").Value = "250d"

# Move selection on ProductLoanInput from B13 to B1
$ws1.Range("B1").Select()

# Move selection/active sheet to ProductLoanOutput (it becomes the visible/active tab)
$ws2.Activate()
$ws2.Range("B1").Select()
